# Applies scheduled-runner market/profit data refresh to Brynhildr_Profits sheets.
$wb = $excel.ActiveWorkbook

# ALC!row17 - One for the Road (Potion)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30336

# ALC!row34 - Sophomore Slump (Goatskin Grimoire)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3296.7144
$ws.Range("I34").Value = 3296.7144
$ws.Range("K34").Value = 3296.7144
$ws.Range("M34").Value = -3093.7144

# ALC!row36 - You Put Your Left Hand In (Engraved Goatskin Grimoire)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 3296.7144
$ws.Range("I36").Value = 3296.7144
$ws.Range("K36").Value = 3296.7144
$ws.Range("M36").Value = -2581.7144

# ALC!row69 - Steeling the Knife, Steeling the Mind (Grade 1 Mind Dissolvent)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8681.666999999999
$ws.Range("J69").Value = 8681.666999999999
$ws.Range("L69").Value = 26045.001
$ws.Range("N69").Value = -27793.001

# ALC!row72 - Surgical Substitution (L) (Grade 1 Mind Dissolvent)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 8681.666999999999
$ws.Range("J72").Value = 8681.666999999999
$ws.Range("L72").Value = 78135.003
$ws.Range("N72").Value = -86871.003

# ALC!row87 - There Was a Late Fee (Noble Gold)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 66666.664
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 66666.664
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 66666.664
$ws.Range("M87").Value = ""
$ws.Range("N87").Value = -69162.664

# ALC!row90 - A Gate Arcane Is Dragon's Bane (L) (Noble Gold)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 66666.664
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 66666.664
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 199999.992
$ws.Range("M90").Value = ""
$ws.Range("N90").Value = -212479.992

# ALC!row107 - Another Man's Ink (Enchanted Truegold Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 13245.375
$ws.Range("I107").Value = 12695.1
$ws.Range("K107").Value = 12695.1
$ws.Range("M107").Value = -10775.1

# ALC!row137 - Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27032618
$ws.Range("I137").Value = 52633020
$ws.Range("J137").Value = 9970.056
$ws.Range("K137").Value = 157899060
$ws.Range("L137").Value = 29910.168
$ws.Range("M137").Value = -157896510
$ws.Range("N137").Value = -35010.16800000001

# ARM!row32 - Ingot We Trust (Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 148918.72
$ws.Range("I32").Value = 186131.52
$ws.Range("J32").Value = 20999.75
$ws.Range("K32").Value = 186131.52
$ws.Range("L32").Value = 20999.75
$ws.Range("M32").Value = -185844.52
$ws.Range("N32").Value = -21573.75

# ARM!row61 - Dealing with the Tough Stuff (Cobalt Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3338187
$ws.Range("I61").Value = 4933.913
$ws.Range("J61").Value = 14290305
$ws.Range("K61").Value = 4933.913
$ws.Range("L61").Value = 14290305
$ws.Range("M61").Value = -4721.913
$ws.Range("N61").Value = -14290729

# ARM!row74 - As the Bolt Flies (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 963093
$ws.Range("I74").Value = 1236688.4
$ws.Range("J74").Value = 16032.077
$ws.Range("K74").Value = 1236688.4
$ws.Range("L74").Value = 16032.077
$ws.Range("M74").Value = -1235814.4
$ws.Range("N74").Value = -17780.077

# ARM!row77 - Heavy Metal Banned (L) (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 963093
$ws.Range("I77").Value = 1236688.4
$ws.Range("J77").Value = 16032.077
$ws.Range("K77").Value = 6183442
$ws.Range("L77").Value = 80160.38499999999
$ws.Range("M77").Value = -6179074
$ws.Range("N77").Value = -88896.38499999999

# ARM!row122 - Haste for High Durium (High Durium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1603.8096
$ws.Range("I122").Value = 1362.4667
$ws.Range("J122").Value = 2207.1667
$ws.Range("K122").Value = 4087.4001
$ws.Range("L122").Value = 6621.500100000001
$ws.Range("M122").Value = -1637.4001
$ws.Range("N122").Value = -11521.5001

# ARM!row136 - Metal with Mettle (Cobalt Tungsten Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3338187
$ws.Range("I136").Value = 4933.913
$ws.Range("J136").Value = 14290305
$ws.Range("K136").Value = 14801.739
$ws.Range("L136").Value = 42870915
$ws.Range("M136").Value = -12251.739
$ws.Range("N136").Value = -42876015

# BSM!row22 - Riveting Run (Iron Rivets)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 282.42856
$ws.Range("I22").Value = 162.83333
$ws.Range("K22").Value = 162.83333
$ws.Range("M22").Value = 10.16667000000001

# CRP!row82 - Aim to Please (Hallowed Chestnut Mask of Aiming)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 42499.5
$ws.Range("J82").Value = 42499.5
$ws.Range("L82").Value = 42499.5
$ws.Range("N82").Value = -43221.5

# CRP!row85 - To Protect My City, I Must Wear a Mask (L) (Hallowed Chestnut Mask of Aiming)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 42499.5
$ws.Range("J85").Value = 42499.5
$ws.Range("L85").Value = 42499.5
$ws.Range("N85").Value = -44995.5

# CUL!row3 - Trout Fishing in Limsa (Grilled Trout)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6670.357
$ws.Range("I3").Value = 2042.7778
$ws.Range("K3").Value = 6128.3334
$ws.Range("M3").Value = -6016.3334

# CUL!row122 - Salt of the North (Northern Sea Salt)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 769048.3
$ws.Range("J122").Value = 1061.909
$ws.Range("L122").Value = 9557.181
$ws.Range("N122").Value = -14457.181

# CUL!row131 - The Mountain Steeped (Tsai tou Vounou)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3769.353
$ws.Range("I131").Value = 1542.5834
$ws.Range("J131").Value = 4983.9546
$ws.Range("K131").Value = 4627.7502
$ws.Range("L131").Value = 14951.8638
$ws.Range("M131").Value = 412.2497999999996
$ws.Range("N131").Value = -25031.8638

# CUL!row138 - Bring Me Your Tacos (Tacos Al Pastor)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 36984.46
$ws.Range("I138").Value = 36984.46
$ws.Range("K138").Value = 110953.38
$ws.Range("M138").Value = -105813.38

# GSM!row132 - On Board for Lar (Lar Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19270
$ws.Range("I132").Value = 17435.889
$ws.Range("J132").Value = 27523.5
$ws.Range("K132").Value = 52307.667
$ws.Range("L132").Value = 82570.5
$ws.Range("M132").Value = -49777.667
$ws.Range("N132").Value = -87630.5

# LTW!row61 - Spelling Me Softly (Raptor Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15656.929
$ws.Range("I61").Value = 14766.417
$ws.Range("J61").Value = 21000
$ws.Range("K61").Value = 14766.417
$ws.Range("L61").Value = 21000
$ws.Range("M61").Value = -14564.417
$ws.Range("N61").Value = -21404

# LTW!row68 - You Could Say It's a Moving Target (Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1952.375
$ws.Range("I68").Value = 1885.7142
$ws.Range("J68").Value = 2419
$ws.Range("K68").Value = 1885.7142
$ws.Range("L68").Value = 2419
$ws.Range("M68").Value = -1136.7142
$ws.Range("N68").Value = -3917

# LTW!row71 - They Call It Bloody Mary (L) (Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1952.375
$ws.Range("I71").Value = 1885.7142
$ws.Range("J71").Value = 2419
$ws.Range("K71").Value = 9428.571
$ws.Range("L71").Value = 12095
$ws.Range("M71").Value = -5684.571
$ws.Range("N71").Value = -19583

# LTW!row93 - Hide to Go Seek (Gagana Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2404.8462
$ws.Range("I93").Value = 1772
$ws.Range("J93").Value = 9999
$ws.Range("K93").Value = 1772
$ws.Range("L93").Value = 9999
$ws.Range("M93").Value = -524
$ws.Range("N93").Value = -12495

# LTW!row113 - Peace in Rest (Atrociraptor Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 15656.929
$ws.Range("I113").Value = 14766.417
$ws.Range("J113").Value = 21000
$ws.Range("K113").Value = 14766.417
$ws.Range("L113").Value = 21000
$ws.Range("M113").Value = -12596.417
$ws.Range("N113").Value = -25340

# LTW!row136 - Respect for Br'aax (Br'aax Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6256304.5
$ws.Range("I136").Value = 4633254
$ws.Range("J136").Value = 9627255
$ws.Range("K136").Value = 13899762
$ws.Range("L136").Value = 28881765
$ws.Range("M136").Value = -13897212
$ws.Range("N136").Value = -28886865

# WVR!row113 - A Tender Table (Pixie Floss)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1166.375
$ws.Range("I113").Value = 1239.6316
$ws.Range("J113").Value = 1059.3077
$ws.Range("K113").Value = 3718.8948
$ws.Range("L113").Value = 3177.9231
$ws.Range("M113").Value = -1548.8948
$ws.Range("N113").Value = -7517.9231

# WVR!row122 - Heavy Armoire (Dark Hempen Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 63479.5
$ws.Range("I122").Value = 1047.9412
$ws.Range("K122").Value = 3143.8236
$ws.Range("M122").Value = -693.8235999999997

# WVR!row132 - Comfy Cabins (Snow Cotton Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4274943.5
$ws.Range("I132").Value = 4630939
$ws.Range("J132").Value = 2997.3333
$ws.Range("K132").Value = 13892817
$ws.Range("L132").Value = 8991.999899999999
$ws.Range("M132").Value = -13890287
$ws.Range("N132").Value = -14051.9999
